# Update "想去人数" (interested-count) values in column F across sheets
# to reflect newly generated numbers (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21426
$ws1.Range("F3").Value = 3359
$ws1.Range("F4").Value = 858
$ws1.Range("F6").Value = 558
$ws1.Range("F7").Value = 808
$ws1.Range("F8").Value = 303
$ws1.Range("F9").Value = 267
$ws1.Range("F10").Value = 74
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 583
$ws1.Range("F13").Value = 192
$ws1.Range("F14").Value = 367
$ws1.Range("F15").Value = 41
$ws1.Range("F16").Value = 466
$ws1.Range("F17").Value = 221
$ws1.Range("F18").Value = 45
$ws1.Range("F19").Value = 32
$ws1.Range("F20").Value = 88
$ws1.Range("F21").Value = 159

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 105

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6179
$ws3.Range("F3").Value = 735
$ws3.Range("F4").Value = 733
$ws3.Range("F5").Value = 1745
$ws3.Range("F6").Value = 92

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6179
$ws4.Range("F3").Value = 735
$ws4.Range("F4").Value = 733
$ws4.Range("F5").Value = 1745
$ws4.Range("F6").Value = 21426
$ws4.Range("F7").Value = 3359
$ws4.Range("F8").Value = 858
$ws4.Range("F10").Value = 92
$ws4.Range("F12").Value = 558
$ws4.Range("F13").Value = 808
$ws4.Range("F14").Value = 303
$ws4.Range("F15").Value = 267
$ws4.Range("F17").Value = 74
$ws4.Range("F20").Value = 144
$ws4.Range("F23").Value = 583
$ws4.Range("F24").Value = 105
$ws4.Range("F25").Value = 192
$ws4.Range("F27").Value = 367
$ws4.Range("F29").Value = 41
$ws4.Range("F30").Value = 466
$ws4.Range("F32").Value = 221
$ws4.Range("F33").Value = 45
$ws4.Range("F36").Value = 32
$ws4.Range("F37").Value = 88
$ws4.Range("F43").Value = 159
